$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.656.64'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.367.01'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +7.54%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.31'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +6.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.27'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +2.75%  '

$ws.Range("E7").Value = '  +23.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.390'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +1.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("B10").Value = 'Cardano'

$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.865'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +10.56%  '

$ws.Range("B11").Value = 'LidoStakedEther'

$ws.Range("C11").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.362.55'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +7.55%  '

$ws.Range("E12").Value = '  -0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '98.424.43'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +1.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.21'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +6.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000248'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +3.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.987.20'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +7.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.48'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.359.69'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +7.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.55'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +0.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.20'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +4.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '487.93'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -5.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.08'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +7.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000211'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +9.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.37'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +6.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.63'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +2.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.91'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.92'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +2.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.548.84'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +7.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.284'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +15.50%  '

$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("E31").Value = '  +8.84%  '

$ws.Range("E32").Value = '  +9.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.64'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +7.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +2.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.00'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +5.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.150'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -1.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.28'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -0.63%  '

$ws.Range("E38").Value = '  +4.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '500.25'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +6.47%  '

$ws.Range("B40").Value = 'WhiteBITCoin'

$ws.Range("C40").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.91'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +2.89%  '

$ws.Range("B41").Value = 'PolygonEcosystemToken'

$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.459'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +5.47%  '

$ws.Range("E42").Value = '  +4.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.26'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +3.54%  '

$ws.Range("E44").Value = '  +5.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.791'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +13.43%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.36'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("E48").Value = '  +1.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.854'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +13.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.61'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +3.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.95'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +3.89%  '
